# Auto-generated Excel COM-interop script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.301.20'
$ws.Range('E2').Value = '  -2.74%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.853.18'
$ws.Range('E3').Value = '  -3.31%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '325.62'
$ws.Range('E5').Value = '  -0.54%  '
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4551'
$ws.Range('E7').Value = '  -2.89%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3895'
$ws.Range('E8').Value = '  -3.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '48.13'
$ws.Range('E9').Value = '  -9.42%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07909'
$ws.Range('E10').Value = '  -5.47%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.010'
$ws.Range('E11').Value = '  -3.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.37'
$ws.Range('E12').Value = '  -3.34%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.850.72'
$ws.Range('E13').Value = '  -2.52%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.907'
$ws.Range('E14').Value = '  -2.50%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.153'
$ws.Range('E15').Value = '  -4.06%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.002'
$ws.Range('E16').Value = '  +0.12%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.06647'
$ws.Range('E17').Value = '  +1.15%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '85.95'
$ws.Range('E18').Value = '  -3.95%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.00001025'
$ws.Range('E19').Value = '  -3.20%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.22'
$ws.Range('E20').Value = '  -3.97%  '
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.489'
$ws.Range('E22').Value = '  -3.88%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.318.72'
$ws.Range('E23').Value = '  -2.62%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.86'
$ws.Range('E24').Value = '  -4.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.283'
$ws.Range('E25').Value = '  +0.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.085.08'
$ws.Range('E26').Value = '  -1.74%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '153.90'
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.93'
$ws.Range('E28').Value = '  -0.24%  '
$ws.Range('E29').Value = '  -3.25%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.451'
$ws.Range('E30').Value = '  -4.28%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '121.08'
$ws.Range('E31').Value = '  -1.56%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9455'
$ws.Range('E32').Value = '  -2.79%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09349'
$ws.Range('E33').Value = '  -2.12%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.440'
$ws.Range('E34').Value = '  -0.36%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.630'
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.251'
$ws.Range('E36').Value = '  -4.96%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06034'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02226'
$ws.Range('E38').Value = '  -3.24%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.214'
$ws.Range('E39').Value = '  -0.12%  '
$ws.Range('B40').Value = 'Frax'
$ws.Range('C40').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.002'
$ws.Range('E40').Value = '  +0.07%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.017'
$ws.Range('E41').Value = '  -8.83%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5924'
$ws.Range('E42').Value = '  -3.14%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1883'
$ws.Range('E43').Value = '  -0.90%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.12'
$ws.Range('E44').Value = '  -8.10%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.278'
$ws.Range('E45').Value = '  -1.86%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5597'
$ws.Range('E46').Value = '  -3.86%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.03'
$ws.Range('E47').Value = '  -5.25%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.377'
$ws.Range('E48').Value = '  -2.02%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.915'
$ws.Range('E49').Value = '  -4.99%  '
$ws.Range('E50').Value = '  -1.48%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '108.13'
$ws.Range('E51').Value = '  -1.09%  '
